# Applies the "coletor" commit: the spotify data-collection run was
# re-executed against a different seed track ("Broken-Hearted Girl" and its
# Alan Braxe remix) instead of the previous Beyoncé titles. The "songs" and
# "credits" sheets are rewritten down to just the rows produced by that run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "songs": keep header + 2 data rows (was header + 16 data rows)
# ---------------------------------------------------------------------
$songs = $wb.Worksheets.Item("songs")

# Drop the now-unused rows 4..17 first so the remaining row count and the
# sheet dimension end up correct.
$songs.Range("A4:N17").EntireRow.Delete()

# Row 2
$songs.Range("A2").Value = 0
$songs.Range("B2").Value = 1
$songs.Range("C2").Value = 278106
$songs.Range("D2").Value = $false
$songs.Range("E2").Value = "{'spotify': 'https://open.spotify.com/track/4iN55SUsXVSh8Og9EDyg3z'}"
$songs.Range("F2").Value = "https://api.spotify.com/v1/tracks/4iN55SUsXVSh8Og9EDyg3z"
$songs.Range("G2").Value = "4iN55SUsXVSh8Og9EDyg3z"
$songs.Range("H2").Value = $false
$songs.Range("I2").Value = $true
$songs.Range("J2").Value = "Broken-Hearted Girl"
$songs.Range("K2").Value = "https://p.scdn.co/mp3-preview/348979c238d82d9faad82ae8064810e2065dae3c?cid=19222f3ec658437489e9280a521db7ad"
$songs.Range("L2").Value = 1
$songs.Range("M2").Value = "track"
$songs.Range("N2").Value = "spotify:track:4iN55SUsXVSh8Og9EDyg3z"

# Row 3
$songs.Range("A3").Value = 1
$songs.Range("B3").Value = 1
$songs.Range("C3").Value = 209240
$songs.Range("D3").Value = $false
$songs.Range("E3").Value = "{'spotify': 'https://open.spotify.com/track/6HtFcQNBIjzOZeNJ0IXKDB'}"
$songs.Range("F3").Value = "https://api.spotify.com/v1/tracks/6HtFcQNBIjzOZeNJ0IXKDB"
$songs.Range("G3").Value = "6HtFcQNBIjzOZeNJ0IXKDB"
$songs.Range("H3").Value = $false
$songs.Range("I3").Value = $true
$songs.Range("J3").Value = "Broken-Hearted Girl - Alan Braxe Remix - Radio Edit"
$songs.Range("K3").Value = "https://p.scdn.co/mp3-preview/788fc83010024e5d2a2405dc2c87c4a3822536c8?cid=19222f3ec658437489e9280a521db7ad"
$songs.Range("L3").Value = 2
$songs.Range("M3").Value = "track"
$songs.Range("N3").Value = "spotify:track:6HtFcQNBIjzOZeNJ0IXKDB"

# ---------------------------------------------------------------------
# Sheet "credits": keep header + 7 data rows (was header + 12 data rows)
# ---------------------------------------------------------------------
$credits = $wb.Worksheets.Item("credits")

# Drop the now-unused rows 9..13 first.
$credits.Range("A9:K13").EntireRow.Delete()

# Row 2 - remixer (J2/K2 were already blank in the source row, leave as-is)
$credits.Range("A2").Value = 0
$credits.Range("B2").Value = "spotify:track:6HtFcQNBIjzOZeNJ0IXKDB"
$credits.Range("C2").Value = "Broken-Hearted Girl - Alan Braxe Remix - Radio Edit"
$credits.Range("D2").Value = "Performers"
$credits.Range("E2").Value = "spotify:artist:24JRvbKfTcF2x7c2kCCJrW"
$credits.Range("F2").Value = "Alan Braxe"
$credits.Range("G2").Value = "https://i.scdn.co/image/ab677762000078e6aae517f57fa4e9833b4bce9c"
$credits.Range("H2").Value = "['remixer']"
$credits.Range("I2").Value = 0.7479000091552734

# Row 3 - main artist (J3/K3 carried a songwriter link before; clear it)
$credits.Range("A3").Value = 1
$credits.Range("B3").Value = "spotify:track:6HtFcQNBIjzOZeNJ0IXKDB"
$credits.Range("C3").Value = "Broken-Hearted Girl - Alan Braxe Remix - Radio Edit"
$credits.Range("D3").Value = "Performers"
$credits.Range("E3").Value = "spotify:artist:6vWDO969PvNqNYHIOW5v0m"
$credits.Range("F3").Value = "Beyoncé"
$credits.Range("G3").Value = "https://i.scdn.co/image/ab677762000078e645c984e8c82f9ce15ebf1f51"
$credits.Range("H3").Value = "['main artist']"
$credits.Range("I3").Value = 0.8999000191688538
$credits.Range("J3").Value = ""
$credits.Range("K3").Value = ""

# Row 4 - writer (Babyface) (J4/K4 carried a songwriter link before; clear it)
$credits.Range("A4").Value = 2
$credits.Range("B4").Value = "spotify:track:6HtFcQNBIjzOZeNJ0IXKDB"
$credits.Range("C4").Value = "Broken-Hearted Girl - Alan Braxe Remix - Radio Edit"
$credits.Range("D4").Value = "Writers"
$credits.Range("E4").Value = "spotify:artist:3aVoqlJOYx31lH1gibGDt3"
$credits.Range("F4").Value = "Babyface"
$credits.Range("G4").Value = "https://i.scdn.co/image/ab677762000078e6d2b377637d9f6ed34f1652e2"
$credits.Range("H4").Value = "['composer', 'lyricist']"
$credits.Range("I4").Value = 0.6998000144958496
$credits.Range("J4").Value = ""
$credits.Range("K4").Value = ""

# Row 5 - writer (Mikkel Storleer Eriksen) (J5/K5 already blank, leave as-is)
$credits.Range("A5").Value = 4
$credits.Range("B5").Value = "spotify:track:6HtFcQNBIjzOZeNJ0IXKDB"
$credits.Range("C5").Value = "Broken-Hearted Girl - Alan Braxe Remix - Radio Edit"
$credits.Range("D5").Value = "Writers"
$credits.Range("E5").Value = "spotify:artist:59mvY7ziUSNoyIGgErHemV"
$credits.Range("F5").Value = "Mikkel Storleer Eriksen"
$credits.Range("G5").Value = "https://i.scdn.co/image/ab677762000078e6afc079cda32d54850e82c385"
$credits.Range("H5").Value = "['composer', 'lyricist']"
$credits.Range("I5").Value = 0.6995999813079834

# Row 6 - writer (Tor Erik Hermansen) (J6/K6 already blank, leave as-is)
$credits.Range("A6").Value = 5
$credits.Range("B6").Value = "spotify:track:6HtFcQNBIjzOZeNJ0IXKDB"
$credits.Range("C6").Value = "Broken-Hearted Girl - Alan Braxe Remix - Radio Edit"
$credits.Range("D6").Value = "Writers"
$credits.Range("E6").Value = "spotify:artist:7BwWmXxaUFHTL8f8IeszOZ"
$credits.Range("F6").Value = "Tor Erik Hermansen"
$credits.Range("G6").Value = "https://i.scdn.co/image/ab677762000078e6afc079cda32d54850e82c385"
$credits.Range("H6").Value = "['composer', 'lyricist']"
$credits.Range("I6").Value = 0.6988999843597412

# Row 7 - producer (StarGate for 45th & 3rd Music LLC) (J7/K7 already blank, leave as-is)
$credits.Range("A7").Value = 6
$credits.Range("B7").Value = "spotify:track:6HtFcQNBIjzOZeNJ0IXKDB"
$credits.Range("C7").Value = "Broken-Hearted Girl - Alan Braxe Remix - Radio Edit"
$credits.Range("D7").Value = "Producers"
$credits.Range("E7").Value = "spotify:artist:0rwyRIAM4qMPtYvkC0pECf"
$credits.Range("F7").Value = "StarGate for 45th & 3rd Music LLC"
$credits.Range("G7").Value = "https://i.scdn.co/image/ab677762000078e6aae517f57fa4e9833b4bce9c"
$credits.Range("H7").Value = "['producer']"
$credits.Range("I7").Value = 0.5491999983787537

# Row 8 - producer (Beyoncé Knowles for Music World Music, LLC) (J8/K8 already blank, leave as-is)
$credits.Range("A8").Value = 7
$credits.Range("B8").Value = "spotify:track:6HtFcQNBIjzOZeNJ0IXKDB"
$credits.Range("C8").Value = "Broken-Hearted Girl - Alan Braxe Remix - Radio Edit"
$credits.Range("D8").Value = "Producers"
$credits.Range("E8").Value = "spotify:artist:0GAqvwpLwRPnWjy5TJKfoT"
$credits.Range("F8").Value = "Beyoncé Knowles for Music World Music, LLC"
$credits.Range("G8").Value = "https://i.scdn.co/image/ab677762000078e6aae517f57fa4e9833b4bce9c"
$credits.Range("H8").Value = "['producer']"
$credits.Range("I8").Value = 0.5490000247955322
